$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear stray forecast values that shouldn't be present
$ws.Range("C2").ClearContents()
$ws.Range("C3").ClearContents()

# Correct minor floating point precision fixes from the naive forecaster bug fix
$ws.Range("C4").Value = -3.956152295564863
$ws.Range("C5").Value = 1.234995474941436
$ws.Range("C6").Value = 0.8993608108207818
$ws.Range("C8").Value = 0.02019328874802717
$ws.Range("E8").Value = -0.03860754389360954
$ws.Range("E10").Value = 0.1987429576382871
$ws.Range("E11").Value = 0.1903092973221998
$ws.Range("C12").Value = 0.0720185131838802
$ws.Range("E13").Value = -0.8612142616933216
$ws.Range("C14").Value = -0.8017595264762423
$ws.Range("E14").Value = 0.06491682578966262
$ws.Range("C16").Value = 0.9704846793491706
$ws.Range("E16").Value = -0.4119900615863981
$ws.Range("C18").Value = 0.3928252664241683
$ws.Range("C19").Value = 0.3224026462283369
$ws.Range("E19").Value = -0.4796777936134977
